$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 31 and 32: two new Mac-Address (usr_id) entries for regcntr_id 10001
$ws.Range("A31").Value = 10001
$ws.Range("B31").Value = 110030
$ws.Range("C31").Value = "eng"
$ws.Range("D31").Value = $true
$ws.Range("E31").Value = "superadmin"
$ws.Range("F31").Value = "now()"
$ws.Range("G31").Value = "now()"

$ws.Range("A32").Value = 10001
$ws.Range("B32").Value = 110031
$ws.Range("C32").Value = "eng"
$ws.Range("D32").Value = $true
$ws.Range("E32").Value = "superadmin"
$ws.Range("F32").Value = "now()"
$ws.Range("G32").Value = "now()"

# Update selection / view to reflect where the user ended up after editing
$ws.Range("E28").Select()
$excel.ActiveWindow.ScrollRow = 19
